# "contingencies with rene fine"
# Two new line rows (line7, line8) are inserted right after line6, pushing
# the existing extr1..extr8 rows down by two. Several C/D/E values are also
# refreshed for the (now shifted) extr rows, as well as for the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing rows 8-15 (extr1..extr8) down to rows 10-17, working
# from the bottom up so rows are not overwritten before being copied, and
# formatting/style indices are preserved via Range.Copy.
for ($r = 15; $r -ge 8; $r--) {
    $target = $r + 2
    $ws.Range("A" + $r + ":E" + $r).Copy($ws.Range("A" + $target + ":E" + $target))
}

# Final target data for rows 2-17: A (index), B (name), C, D, E (in_service)
$data = @(
    @(2,  0,  "line1", 7,  9,  $true),
    @(3,  1,  "line2", 9,  8,  $true),
    @(4,  2,  "line3", 8,  10, $true),
    @(5,  3,  "line4", 8,  11, $false),
    @(6,  4,  "line5", 10, 5,  $true),
    @(7,  5,  "line6", 12, 8,  $true),
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $false),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
    $ws.Cells.Item($r, 5).Value = $item[5]
}
